# Add the "2022-Q3" worksheet (new quarterly snapshot) right after the
# "总计" (summary) sheet, and update the summary sheet with the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right after "总计" (i.e. as the 2nd tab)
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $firstSheet)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Columns that hold numeric-looking text (fund code / scale / position /
# ratio / market value) must be forced to Text format so Excel does not
# coerce them into numbers (which would lose leading zeros / formatting).
$q3.Range("B2:B11").NumberFormat = "@"
$q3.Range("D2:G11").NumberFormat = "@"

$q3Data = @(
    @(0, "206009", "鹏华新兴产业混合",           "41.17", "89.21", "9.79", "4.0305", 1),
    @(1, "012930", "中庚价值先锋股票",           "47.83", "94.71", "3.91", "1.8702", 9),
    @(2, "001468", "广发改革先锋灵活配置混合",   "6.25",  "62.12", "2.17", "0.1356", 10),
    @(3, "920002", "中金精选股票A",              "2.95",  "82.28", "2.44", "0.0720", 9),
    @(4, "014126", "华夏中证1000指数增强C",      "8.78",  "89.62", "0.80", "0.0702", 6),
    @(5, "350002", "天治低碳经济灵活配置混合",   "0.77",  "63.67", "1.93", "0.0149", 7),
    @(6, "014125", "华夏中证1000指数增强A",      "0.97",  "89.62", "0.80", "0.0078", 6),
    @(7, "920922", "中金精选股票C",              "0.11",  "82.28", "2.44", "0.0027", 9),
    @(8, "006143", "恒生前海中证质量成长低波动指数A", "0.05", "93.33", "2.48", "0.0012", 5),
    @(9, "006144", "恒生前海中证质量成长低波动指数C", "0.01", "93.33", "2.48", "0.0002", 5)
)

$row = 2
foreach ($rec in $q3Data) {
    $q3.Range("A$row").Value = $rec[0]
    $q3.Range("B$row").Value = $rec[1]
    $q3.Range("C$row").Value = $rec[2]
    $q3.Range("D$row").Value = $rec[3]
    $q3.Range("E$row").Value = $rec[4]
    $q3.Range("F$row").Value = $rec[5]
    $q3.Range("G$row").Value = $rec[6]
    $q3.Range("H$row").Value = $rec[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q3 on
#    top of the existing data (pushing all older quarters down by one).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The Insert() call can leave stray formatting behind on the new blank
# row; clear it before writing fresh values.
$summary.Range("A2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 10
$summary.Range("D2").Value = 6.21

# Match the styling used by the other index cells in column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Keep the last tab ("2021-Q2") as the active / selected sheet, same
#    as before the edit.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
